$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.420.29'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '3.139.54'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'607.79"
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").Value = "'143.95"
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D8").Value = '3.134.49'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").Value = "'5.38"
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("E13").Value = '  +3.91%  '
$ws.Range("D14").Value = "'35.40"
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '3.658.22'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").Value = '64.385.28'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").Value = '3.142.27'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = "'6.85"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = "'477.41"
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").Value = "'14.77"
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = "'0.716"
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("D23").Value = "'7.72"
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").Value = "'85.64"
$ws.Range("E24").Value = '  +3.39%  '
$ws.Range("D25").Value = "'13.40"
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -3.57%  '
$ws.Range("D28").Value = "'8.45"
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = '  +8.21%  '
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("D31").Value = "'2.05"
$ws.Range("E31").Value = '  -5.52%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").Value = "'26.89"
$ws.Range("E33").Value = '  +3.27%  '
$ws.Range("D34").Value = "'2.63"
$ws.Range("E34").Value = '  -2.93%  '
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("D37").Value = '0.0₃0766'
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("D38").Value = "'52.46"
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").Value = "'3.03"
$ws.Range("E39").Value = '  +3.96%  '
$ws.Range("D40").Value = "'445.70"
$ws.Range("E40").Value = '  -2.73%  '
$ws.Range("D41").Value = "'0.0392"
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").Value = "'0.121"
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("D43").Value = "'8.25"
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("D44").Value = '2.883.66'
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("D45").Value = "'0.261"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("E47").Value = '  +3.51%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "'26.20"
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = "'33.73"
$ws.Range("E51").Value = '  +8.30%  '
